$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.137.84"
$ws.Range("D3").Value = "3.531.31"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.79"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.97"
$ws.Range("D7").Value = "3.530.34"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D13").Value = "4.136.11"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.22"
$ws.Range("D16").Value = "3.532.00"
$ws.Range("D18").Value = "65.237.18"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.28"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.96"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.574"
$ws.Range("D24").Value = "3.677.63"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.74"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.14"
$ws.Range("D32").Value = "3.545.42"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.82"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.25"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.99"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.823"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.07"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("D50").Value = "2.381.79"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "303.92"

$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +3.17%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").Value = "  +4.82%  "
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("E21").Value = "  +5.07%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +3.61%  "
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +7.82%  "
$ws.Range("E28").Value = "  +8.24%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("E36").Value = "  +9.47%  "
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  +16.79%  "
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("E48").Value = "  +6.50%  "
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("E50").Value = "  +8.62%  "
$ws.Range("E51").Value = "  +7.06%  "

Write-Output "Updated cryptos list"